$wb = $excel.ActiveWorkbook

# --- MassInertia sheet: remove the extra "lltd" setup rows (3-5), keep only the baseline row ---
$wsMass = $wb.Worksheets.Item("MassInertia")
[void]$wsMass.Rows("3:5").Delete()
[void]$wsMass.Range("H2").Select()

# --- Aero sheet: bump frontal_area-ish column B value from 3.5 to 5 ---
$wsAero = $wb.Worksheets.Item("Aero")
$wsAero.Range("B2").Value = 5
[void]$wsAero.Range("B3").Select()

# --- Engine sheet: final_gear_ratio now computed (1.1*2.9) and maximum_power raised to 550 ---
$wsEngine = $wb.Worksheets.Item("Engine")
$wsEngine.Range("A2:A5").Formula = "=1.1*2.9"
$wsEngine.Range("C2:C5").Value = 550
[void]$wsEngine.Range("D9").Select()

# --- Brake sheet: raise maximum_torque values to a flat 6000 across all setups ---
$wsBrake = $wb.Worksheets.Item("Brake")
$wsBrake.Range("B2:B5").Value = 6000
[void]$wsBrake.Range("E6").Select()

# Brake becomes the active/visible tab when the workbook is reopened.
[void]$wsBrake.Activate()
